$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.930.06"
$ws.Range("E2").Value = "  +3.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.539.53"
$ws.Range("E3").Value = "  +3.33%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.26"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.93"
$ws.Range("E6").Value = "  +4.13%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.538.99"
$ws.Range("E9").Value = "  +3.26%  "
$ws.Range("E10").Value = "  +2.12%  "
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.30"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.30"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("E15").Value = "  +3.17%  "
$ws.Range("E16").Value = "  +3.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.680.61"
$ws.Range("E17").Value = "  +3.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.537.29"
$ws.Range("E18").Value = "  +2.85%  "
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("E21").Value = "  +3.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "330.21"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("E23").Value = "  +3.02%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.19"
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "643.74"
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("E28").Value = "  +8.72%  "
$ws.Range("E30").Value = "  +5.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("E32").Value = "  +2.43%  "
$ws.Range("E33").Value = "  +3.01%  "
$ws.Range("E34").Value = "  +4.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  +4.42%  "
$ws.Range("E37").Value = "  +2.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.64"
$ws.Range("E38").Value = "  +6.59%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.87"
$ws.Range("E39").Value = "  +5.51%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.99"
$ws.Range("E40").Value = "  +3.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.372"
$ws.Range("E41").Value = "  +1.55%  "
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("E43").Value = "  +6.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "161.83"
$ws.Range("E44").Value = "  +5.98%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.64"
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("E48").Value = "  +2.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.60"
$ws.Range("E49").Value = "  +6.82%  "
$ws.Range("E50").Value = "  +3.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0519"
$ws.Range("E51").Value = "  +2.66%  "
